$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9677374958992004
$ws.Range("B1").Value = 0.904264509677887
$ws.Range("C1").Value = 0.7815846800804138
$ws.Range("D1").Value = 0.83417809009552
$ws.Range("E1").Value = 0.9995025396347046
